# Add pipe and colon
# `|` and `:` broke into rdf, manually escape them
#
# Populate the new escaping test rows (16-32) in the exact order that
# reproduces the recorded shared-string table order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 'exclamation'
$ws.Range("A17").Value = 'question'
$ws.Range("B17").Value = 'v?v'
$ws.Range("B16").Value = 'v!v'
$ws.Range("A18").Value = 'doubleQuestion'
$ws.Range("B18").Value = 'v?v?v'
$ws.Range("B19").Value = 'v+v'
$ws.Range("B20").Value = 'v-v'
$ws.Range("A19").Value = 'plus'
$ws.Range("A20").Value = 'minus'
$ws.Range("A21").Value = 'star'
$ws.Range("B21").Value = 'v*v'
$ws.Range("A22").Value = 'hash'
$ws.Range("B22").Value = 'v#v'
$ws.Range("A23").Value = 'doubleHash'
$ws.Range("B23").Value = 'v#v#v'
$ws.Range("A24").Value = 'percent'
$ws.Range("B24").Value = 'v%v'
$ws.Range("A25").Value = 'ampersand'
$ws.Range("B25").Value = 'v&v'
$ws.Range("A26").Value = 'equal'
$ws.Range("B26").Value = 'v=v'
$ws.Range("B27").Value = 'v@v'
$ws.Range("A28").Value = 'colon'
$ws.Range("B28").Value = 'v:v'
$ws.Range("A27").Value = 'at'
$ws.Range("A29").Value = 'gt'
$ws.Range("A30").Value = 'lt'
$ws.Range("A31").Value = 'pipe'
$ws.Range("A32").Value = 'semicolon'
$ws.Range("B29").Value = 'v>v'
$ws.Range("B31").Value = 'v|v'
$ws.Range("B32").Value = 'v;v'
$ws.Range("B30").Value = 'v<v'

# `v@v` (the "at" row) is an email-like token, so Excel auto-hyperlinked it
# with a mailto: link, which also applies the built-in Hyperlink cell style.
$ws.Hyperlinks.Add($ws.Range("B27"), "mailto:v@v") | Out-Null

# Selection ends on D33 after the last entry.
$ws.Range("D33").Select() | Out-Null
